# Lab Exam 03 grading workbook: fill in the "Total Points" (column E)
# scores that mirror the "Points for grading" (column D) values for the
# Customer Class and Product Class rubric sections, and move the cursor
# to the reviewer's final resting cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Customer Class section (rows 3-6): award full points, matching column D.
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Product Class section (rows 10-14): award full points, matching column D.
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Leave the cursor on E15 (the Product Class total), matching the
# reviewer's final selection, and scroll the view back to the top.
$ws.Activate()
$ws.Range("E15").Select() | Out-Null
